$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '45.764.23'
$ws.Range('E2').Value = '  +6.41%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.413.54'
$ws.Range('E3').Value = '  +4.84%  '

$ws.Range('E4').Value = '  +0.21%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '115.60'
$ws.Range('E5').Value = '  +10.14%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '319.21'
$ws.Range('E6').Value = '  +2.20%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.639'
$ws.Range('E7').Value = '  +2.43%  '

$ws.Range('E9').Value = '  +4.49%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '43.06'
$ws.Range('E10').Value = '  +8.16%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0934'
$ws.Range('E11').Value = '  +3.67%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.79'
$ws.Range('E12').Value = '  +6.27%  '

$ws.Range('E13').Value = '  +2.88%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.02'
$ws.Range('E14').Value = '  +3.50%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.98'
$ws.Range('E15').Value = '  +4.17%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.781.49'
$ws.Range('E16').Value = '  -1.78%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.457.19'
$ws.Range('E17').Value = '  +7.51%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '45.791.12'
$ws.Range('E18').Value = '  +7.02%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.60'
$ws.Range('E19').Value = '  +4.07%  '

$ws.Range('E20').Value = '  +4.36%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.41'
$ws.Range('E21').Value = '  -1.52%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '75.22'
$ws.Range('E22').Value = '  +2.44%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.56'
$ws.Range('E23').Value = '  +2.50%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '267.81'
$ws.Range('E24').Value = '  -0.21%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.41'
$ws.Range('E25').Value = '  +10.06%  '

$ws.Range('E26').Value = '  -0.72%  '

$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.41'
$ws.Range('E27').Value = '  +5.33%  '

$ws.Range('B28').Value = 'Filecoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.62'
$ws.Range('E28').Value = '  +6.40%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '40.27'
$ws.Range('E30').Value = '  +11.65%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.97'
$ws.Range('E31').Value = '  +2.73%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0983'
$ws.Range('E32').Value = '  +15.32%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '173.41'
$ws.Range('E33').Value = '  +5.31%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.97'
$ws.Range('E34').Value = '  +13.13%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.07'
$ws.Range('E35').Value = '  +11.20%  '

$ws.Range('E36').Value = '  +2.38%  '

$ws.Range('E37').Value = '  +7.73%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.25'
$ws.Range('E38').Value = '  +16.42%  '

$ws.Range('E39').Value = '  +12.24%  '

$ws.Range('E40').Value = '  +5.88%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.81'
$ws.Range('E41').Value = '  +14.04%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '101.12'
$ws.Range('E42').Value = '  -7.13%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '13.70'
$ws.Range('E43').Value = '  +12.48%  '

$ws.Range('B44').Value = 'MultiversX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '73.09'
$ws.Range('E44').Value = '  +3.08%  '

$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.241'
$ws.Range('E45').Value = '  +6.24%  '

$ws.Range('E46').Value = '  -0.40%  '

$ws.Range('B47').Value = 'ordi'
$ws.Range('C47').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '85.72'
$ws.Range('E47').Value = '  +9.80%  '

$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.84'
$ws.Range('E48').Value = '  +13.48%  '

$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '117.55'
$ws.Range('E49').Value = '  +6.27%  '

$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.54'
$ws.Range('E50').Value = '  +10.28%  '

$ws.Range('B51').Value = 'MinaProtocolToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.66'
$ws.Range('E51').Value = '  +17.53%  '
